# Insert a new "pt_max" column between pt_min (D) and boson (E), pushing
# every subsequent column one slot to the right (E->F, F->G, ... M->N).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at E; Excel shifts E:M -> F:N automatically.
$ws.Columns("E:E").Insert()

# New column header + data (same value, 50, for every data row).
$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E12").Value = 50

# Match the author's final cell selection.
[void]$ws.Range("E17").Select()
